$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1072.25
$ws.Range("I96").Value = 1447.6
$ws.Range("J96").Value = 696.9
$ws.Range("K96").Value = 4342.799999999999
$ws.Range("L96").Value = 2090.7
$ws.Range("M96").Value = -2969.799999999999
$ws.Range("N96").Value = -4836.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13788.815
$ws.Range("I32").Value = 9477.014999999999
$ws.Range("J32").Value = 45887.777
$ws.Range("K32").Value = 9477.014999999999
$ws.Range("L32").Value = 45887.777
$ws.Range("M32").Value = -9190.014999999999
$ws.Range("N32").Value = -46461.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2323.625
$ws.Range("I61").Value = 2047.3636
$ws.Range("J61").Value = 2931.4
$ws.Range("K61").Value = 2047.3636
$ws.Range("L61").Value = 2931.4
$ws.Range("M61").Value = -1835.3636
$ws.Range("N61").Value = -3355.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 678.5625
$ws.Range("I97").Value = 557.13336
$ws.Range("K97").Value = 557.13336
$ws.Range("M97").Value = -61.13336000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1882.7273
$ws.Range("I102").Value = 1744.9445
$ws.Range("J102").Value = 2502.75
$ws.Range("K102").Value = 1744.9445
$ws.Range("L102").Value = 2502.75
$ws.Range("M102").Value = -122.9445000000001
$ws.Range("N102").Value = -5746.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2323.625
$ws.Range("I136").Value = 2047.3636
$ws.Range("J136").Value = 2931.4
$ws.Range("K136").Value = 6142.0908
$ws.Range("L136").Value = 8794.200000000001
$ws.Range("M136").Value = -3592.0908
$ws.Range("N136").Value = -13894.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1673.6957
$ws.Range("I86").Value = 1400.375
$ws.Range("J86").Value = 2298.4285
$ws.Range("K86").Value = 1400.375
$ws.Range("L86").Value = 2298.4285
$ws.Range("M86").Value = -277.375
$ws.Range("N86").Value = -4544.4285

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1673.6957
$ws.Range("I89").Value = 1400.375
$ws.Range("J89").Value = 2298.4285
$ws.Range("K89").Value = 7001.875
$ws.Range("L89").Value = 11492.1425
$ws.Range("M89").Value = -1385.875
$ws.Range("N89").Value = -22724.1425

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 40001064
$ws.Range("I94").Value = 590.9286
$ws.Range("J94").Value = 90910760
$ws.Range("K94").Value = 590.9286
$ws.Range("L94").Value = 90910760
$ws.Range("M94").Value = -139.9286
$ws.Range("N94").Value = -90911662

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1805
$ws.Range("I99").Value = 1539.4445
$ws.Range("K99").Value = 1539.4445
$ws.Range("M99").Value = -41.44450000000006

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1625570.1
$ws.Range("I105").Value = 2842672.8
$ws.Range("J105").Value = 2766.6667
$ws.Range("K105").Value = 2842672.8
$ws.Range("L105").Value = 2766.6667
$ws.Range("M105").Value = -2840925.8
$ws.Range("N105").Value = -6260.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 52635004
$ws.Range("I105").Value = 90912456
$ws.Range("J105").Value = 3512.5
$ws.Range("K105").Value = 90912456
$ws.Range("L105").Value = 3512.5
$ws.Range("M105").Value = -90910709
$ws.Range("N105").Value = -7006.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 10639421
$ws.Range("I131").Value = 100000260
$ws.Range("J131").Value = 1226.4048
$ws.Range("K131").Value = 300000780
$ws.Range("L131").Value = 3679.2144
$ws.Range("M131").Value = -299995740
$ws.Range("N131").Value = -13759.2144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2715.7896
$ws.Range("I140").Value = 1393.0769
$ws.Range("J140").Value = 5581.6665
$ws.Range("K140").Value = 4179.2307
$ws.Range("L140").Value = 16744.9995
$ws.Range("M140").Value = 1000.7693
$ws.Range("N140").Value = -27104.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2929.8
$ws.Range("I80").Value = 2206.25
$ws.Range("J80").Value = 3270.2942
$ws.Range("K80").Value = 2206.25
$ws.Range("L80").Value = 3270.2942
$ws.Range("M80").Value = -1208.25
$ws.Range("N80").Value = -5266.2942

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2929.8
$ws.Range("I83").Value = 2206.25
$ws.Range("J83").Value = 3270.2942
$ws.Range("K83").Value = 11031.25
$ws.Range("L83").Value = 16351.471
$ws.Range("M83").Value = -6039.25
$ws.Range("N83").Value = -26335.471

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1010
$ws.Range("I97").Value = 719.0909
$ws.Range("J97").Value = 1467.1428
$ws.Range("K97").Value = 719.0909
$ws.Range("L97").Value = 1467.1428
$ws.Range("M97").Value = -223.0909
$ws.Range("N97").Value = -2459.1428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 25600
$ws.Range("I137").Value = 10000
$ws.Range("J137").Value = 33400
$ws.Range("K137").Value = 10000
$ws.Range("L137").Value = 33400
$ws.Range("M137").Value = -4900
$ws.Range("N137").Value = -43600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 81133.336
$ws.Range("J138").Value = 81133.336
$ws.Range("L138").Value = 81133.336
$ws.Range("N138").Value = -91413.336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 59072
$ws.Range("J140").Value = 58840
$ws.Range("L140").Value = 58840
$ws.Range("N140").Value = -69200

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 16926
$ws.Range("I68").Value = 34600.668
$ws.Range("J68").Value = 3670
$ws.Range("K68").Value = 34600.668
$ws.Range("L68").Value = 3670
$ws.Range("M68").Value = -33851.668
$ws.Range("N68").Value = -5168

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 16926
$ws.Range("I71").Value = 34600.668
$ws.Range("J71").Value = 3670
$ws.Range("K71").Value = 173003.34
$ws.Range("L71").Value = 18350
$ws.Range("M71").Value = -169259.34
$ws.Range("N71").Value = -25838

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2399.6843
$ws.Range("I82").Value = 1778.8
$ws.Range("J82").Value = 2621.4285
$ws.Range("K82").Value = 1778.8
$ws.Range("L82").Value = 2621.4285
$ws.Range("M82").Value = -1417.8
$ws.Range("N82").Value = -3343.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2399.6843
$ws.Range("I85").Value = 1778.8
$ws.Range("J85").Value = 2621.4285
$ws.Range("K85").Value = 1778.8
$ws.Range("L85").Value = 2621.4285
$ws.Range("M85").Value = -530.8
$ws.Range("N85").Value = -5117.4285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5122.533
$ws.Range("I93").Value = 6256
$ws.Range("J93").Value = 2477.7778
$ws.Range("K93").Value = 6256
$ws.Range("L93").Value = 2477.7778
$ws.Range("M93").Value = -5008
$ws.Range("N93").Value = -4973.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2633.25
$ws.Range("I100").Value = 2599.9092
$ws.Range("K100").Value = 2599.9092
$ws.Range("M100").Value = -2058.9092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4128.5713
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4150
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4150
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5398

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4128.5713
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4150
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 20750
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -26990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8963.5
$ws.Range("I81").Value = 100000
$ws.Range("J81").Value = 1960.6923
$ws.Range("K81").Value = 200000
$ws.Range("L81").Value = 3921.3846
$ws.Range("M81").Value = -198939
$ws.Range("N81").Value = -6043.384599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8963.5
$ws.Range("I84").Value = 100000
$ws.Range("J84").Value = 1960.6923
$ws.Range("K84").Value = 1000000
$ws.Range("L84").Value = 19606.923
$ws.Range("M84").Value = -994696
$ws.Range("N84").Value = -30214.923

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1914.4
$ws.Range("I96").Value = 1756
$ws.Range("J96").Value = 2350
$ws.Range("K96").Value = 1756
$ws.Range("L96").Value = 2350
$ws.Range("M96").Value = -383
$ws.Range("N96").Value = -5096
